# Apply updated dSF (column F) values for specific rows.
# These reflect a data repull / recalculation of the mean for the dSF column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    17 = 1
    22 = -1
    31 = 1
    34 = 2
    36 = -1
    42 = -3
    45 = -3
    47 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
